$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet updates ---
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("A2").Value = "Version: $newVersion"

$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Huashan Coal Mine, China, M1965, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$aboutSheet.Range("A6").Value = $newCitation

# --- "Boundaries and methane sources" sheet updates ---
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 12; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
